# Refresh Universalis market-price snapshot values across the Leve profit sheets
# (mirrors an automated "scheduled runner" data sync commit).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 497.8
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 422.25
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 422.25
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -1074.25
$ws.Range("H112").Value = 1394.29
$ws.Range("J112").Value = 1394.2424
$ws.Range("L112").Value = 4182.7272
$ws.Range("N112").Value = -6398.7272
$ws.Range("H115").Value = 1502.3077
$ws.Range("I115").Value = 1420.909
$ws.Range("J115").Value = 1950
$ws.Range("K115").Value = 4262.727000000001
$ws.Range("L115").Value = 5850
$ws.Range("M115").Value = -2695.727000000001
$ws.Range("N115").Value = -8984
$ws.Range("H118").Value = 258
$ws.Range("I118").Value = 258
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 774
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 883
$ws.Range("H129").Value = 874.5833
$ws.Range("J129").Value = 963.0244
$ws.Range("L129").Value = 2889.0732
$ws.Range("N129").Value = -12889.0732
$ws.Range("N118").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 5752.5
$ws.Range("I110").Value = 7170
$ws.Range("K110").Value = 7170
$ws.Range("M110").Value = -5125
$ws.Range("H122").Value = 2341.6
$ws.Range("I122").Value = 1445.6471
$ws.Range("K122").Value = 4336.9413
$ws.Range("M122").Value = -1886.9413

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 28000
$ws.Range("I69").Value = 15000
$ws.Range("J69").Value = 47500
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 47500
$ws.Range("M69").Value = -14251
$ws.Range("N69").Value = -48998
$ws.Range("H72").Value = 28000
$ws.Range("I72").Value = 15000
$ws.Range("J72").Value = 47500
$ws.Range("K72").Value = 45000
$ws.Range("L72").Value = 142500
$ws.Range("M72").Value = -41256
$ws.Range("N72").Value = -149988

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 8901.733
$ws.Range("I34").Value = 322
$ws.Range("J34").Value = 12021.637
$ws.Range("K34").Value = 966
$ws.Range("L34").Value = 36064.911
$ws.Range("M34").Value = -882
$ws.Range("N34").Value = -36232.911
$ws.Range("H39").Value = 7927.8
$ws.Range("J39").Value = 8294
$ws.Range("L39").Value = 24882
$ws.Range("N39").Value = -25470
$ws.Range("H112").Value = 7689.5713
$ws.Range("J112").Value = 10000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32216
$ws.Range("H121").Value = 1911.7424
$ws.Range("I121").Value = 289.85715
$ws.Range("J121").Value = 2104.1694
$ws.Range("K121").Value = 869.5714499999999
$ws.Range("L121").Value = 6312.5082
$ws.Range("M121").Value = 440.4285500000001
$ws.Range("N121").Value = -8932.5082
$ws.Range("H129").Value = 3406
$ws.Range("I129").Value = 6176.6665
$ws.Range("J129").Value = 1558.8889
$ws.Range("K129").Value = 18529.9995
$ws.Range("L129").Value = 4676.6667
$ws.Range("M129").Value = -13529.9995
$ws.Range("N129").Value = -14676.6667
$ws.Range("H134").Value = 3907.96
$ws.Range("I134").Value = 3007.6155
$ws.Range("K134").Value = 9022.8465
$ws.Range("M134").Value = -3952.8465
$ws.Range("H136").Value = 2877.5
$ws.Range("I136").Value = 2591.25
$ws.Range("J136").Value = 3450
$ws.Range("K136").Value = 7773.75
$ws.Range("L136").Value = 10350
$ws.Range("M136").Value = -2673.75
$ws.Range("N136").Value = -20550
$ws.Range("H139").Value = 3964.5
$ws.Range("I139").Value = 1943.5
$ws.Range("K139").Value = 5830.5
$ws.Range("M139").Value = -690.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6153.4614
$ws.Range("I7").Value = 4116.6665
$ws.Range("J7").Value = 7899.2856
$ws.Range("K7").Value = 4116.6665
$ws.Range("L7").Value = 7899.2856
$ws.Range("M7").Value = -4004.6665
$ws.Range("N7").Value = -8123.2856
$ws.Range("H22").Value = 3005.9412
$ws.Range("I22").Value = 3133.3333
$ws.Range("J22").Value = 2936.4546
$ws.Range("K22").Value = 3133.3333
$ws.Range("L22").Value = 2936.4546
$ws.Range("M22").Value = -2838.3333
$ws.Range("N22").Value = -3526.4546
$ws.Range("H27").Value = 3005.9412
$ws.Range("I27").Value = 3133.3333
$ws.Range("J27").Value = 2936.4546
$ws.Range("K27").Value = 3133.3333
$ws.Range("L27").Value = 2936.4546
$ws.Range("M27").Value = -3026.3333
$ws.Range("N27").Value = -3150.4546
$ws.Range("H40").Value = 6104.467
$ws.Range("I40").Value = 5723.0835
$ws.Range("J40").Value = 7630
$ws.Range("K40").Value = 5723.0835
$ws.Range("L40").Value = 7630
$ws.Range("M40").Value = -5587.0835
$ws.Range("N40").Value = -7902
$ws.Range("H46").Value = 2165.3044
$ws.Range("I46").Value = 2183.3333
$ws.Range("J46").Value = 2158.9412
$ws.Range("K46").Value = 2183.3333
$ws.Range("L46").Value = 2158.9412
$ws.Range("M46").Value = -1995.3333
$ws.Range("N46").Value = -2534.9412
$ws.Range("H54").Value = 30083.5
$ws.Range("J54").Value = 30083.5
$ws.Range("L54").Value = 30083.5
$ws.Range("N54").Value = -31371.5
$ws.Range("H63").Value = 45257
$ws.Range("J63").Value = 45257
$ws.Range("L63").Value = 45257
$ws.Range("N63").Value = -46755
$ws.Range("H66").Value = 45257
$ws.Range("J66").Value = 45257
$ws.Range("L66").Value = 135771
$ws.Range("N66").Value = -143259
$ws.Range("H68").Value = 896.4761999999999
$ws.Range("I68").Value = 731.7027
$ws.Range("J68").Value = 2115.8
$ws.Range("K68").Value = 731.7027
$ws.Range("L68").Value = 2115.8
$ws.Range("M68").Value = 17.29729999999995
$ws.Range("N68").Value = -3613.8
$ws.Range("H71").Value = 896.4761999999999
$ws.Range("I71").Value = 731.7027
$ws.Range("J71").Value = 2115.8
$ws.Range("K71").Value = 3658.5135
$ws.Range("L71").Value = 10579
$ws.Range("M71").Value = 85.48649999999998
$ws.Range("N71").Value = -18067
$ws.Range("H74").Value = 44000
$ws.Range("J74").Value = 44000
$ws.Range("L74").Value = 44000
$ws.Range("N74").Value = -45996
$ws.Range("H77").Value = 44000
$ws.Range("J77").Value = 44000
$ws.Range("L77").Value = 132000
$ws.Range("N77").Value = -141984
$ws.Range("H123").Value = 29979.143
$ws.Range("J123").Value = 29979.143
$ws.Range("L123").Value = 29979.143
$ws.Range("N123").Value = -39779.143
$ws.Range("H126").Value = 6153.4614
$ws.Range("I126").Value = 4116.6665
$ws.Range("J126").Value = 7899.2856
$ws.Range("K126").Value = 12349.9995
$ws.Range("L126").Value = 23697.8568
$ws.Range("M126").Value = -9879.999500000002
$ws.Range("N126").Value = -28637.8568

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 771708.0600000001
$ws.Range("I126").Value = 4215.6
$ws.Range("J126").Value = 1198092.8
$ws.Range("K126").Value = 12646.8
$ws.Range("L126").Value = 3594278.4
$ws.Range("M126").Value = -10176.8
$ws.Range("N126").Value = -3599218.4
